# Add the new DIF row for the 2020-12-31 report to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 4: Date, SecurityId, Liquidity
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats, carry over the date style used by A2/A3
$ws.Range("A4").Value = 44196
$ws.Range("B4").Value = "6618 HK Equity"
$ws.Range("C4").Value = "L0"

# Move the selection the way Excel left it after the edit.
$ws.Range("A5").Select()

$wb.Save()
